# Update data: 2025-11-18 09:20
# "Stock List" sheet gains two new rows at the top (GROWW, TMCV); every
# existing row shifts down by two, and the two rows that fall off the
# bottom (old DIGITIDE / SCODATUBES rows) are dropped so the sheet keeps
# its original 76-row extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock List")

# Make room for the two new rows right after the header.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# The freshly inserted rows inherit the header's bold/centered style from
# the insert operation - strip that back to the plain data-row look.
$ws.Range("A2:H3").ClearFormats()

$ws.Range("A2").Value = "📋"
$ws.Range("B2").Value = "GROWW"
$ws.Range("C2").Value = "GROWW"
$ws.Range("D2").Value = 185
$ws.Range("E2").Value = 6.0476
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = 107698.3932

$ws.Range("A3").Value = "📋"
$ws.Range("B3").Value = "TMCV"
$ws.Range("C3").Value = "TMCV"
$ws.Range("D3").Value = 317.2
$ws.Range("E3").Value = -1.6739
$ws.Range("F3").Value = "N/A"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = 118792.0101

# Drop the two rows that overflowed past the bottom of the table.
$ws.Rows.Item(78).Delete()
$ws.Rows.Item(77).Delete()
